# Fruta / hortaliza, semanal
# Insert two new weekly rows for "Membrillo" (row 62 and 63), shifting
# all subsequent rows down by two (old row 62 -> new row 64, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 62 (shifts 62.. down to 64..)
$ws.Rows("62:63").Insert()

# New row 62 data
$ws.Cells.Item(62, 1).Value = 6
$ws.Cells.Item(62, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(62, 3).Value = "Metropolitana"
$ws.Cells.Item(62, 4).Value = 44673
$ws.Cells.Item(62, 5).Value = 13
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100104
$ws.Cells.Item(62, 8).Value = "Frutos de pepita"
$ws.Cells.Item(62, 9).Value = 100104003
$ws.Cells.Item(62, 10).Value = "Membrillo"
$ws.Cells.Item(62, 11).Value = "Champion"
$ws.Cells.Item(62, 12).Value = "Primera"
$ws.Cells.Item(62, 13).Value = 20
$ws.Cells.Item(62, 14).Value = 222000
$ws.Cells.Item(62, 15).Value = 230000
$ws.Cells.Item(62, 16).Value = 226000
$ws.Cells.Item(62, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(62, 18).Value = "Provincia de Cachapoal"
$ws.Cells.Item(62, 19).Value = 565
$ws.Cells.Item(62, 20).Value = 400

# New row 63 data
$ws.Cells.Item(63, 1).Value = 6
$ws.Cells.Item(63, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(63, 3).Value = "Metropolitana"
$ws.Cells.Item(63, 4).Value = 44673
$ws.Cells.Item(63, 5).Value = 13
$ws.Cells.Item(63, 6).Value = "Fruta"
$ws.Cells.Item(63, 7).Value = 100104
$ws.Cells.Item(63, 8).Value = "Frutos de pepita"
$ws.Cells.Item(63, 9).Value = 100104003
$ws.Cells.Item(63, 10).Value = "Membrillo"
$ws.Cells.Item(63, 11).Value = "Champion"
$ws.Cells.Item(63, 12).Value = "Primera"
$ws.Cells.Item(63, 13).Value = 8
$ws.Cells.Item(63, 14).Value = 220000
$ws.Cells.Item(63, 15).Value = 220000
$ws.Cells.Item(63, 16).Value = 220000
$ws.Cells.Item(63, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(63, 18).Value = "Región Metropolitana"
$ws.Cells.Item(63, 19).Value = 489
$ws.Cells.Item(63, 20).Value = 450
